$d = $word.ActiveDocument

# --- Change 1: mark the figure's run as noProof (adds <w:rPr><w:noProof/></w:rPr>
#     to the run that holds <w:lastRenderedPageBreak/><w:drawing>...) ---
$shape = $d.InlineShapes.Item(1)
$shape.Range.NoProofing = $true

# --- Change 2: fix the "llitbang" typo -> "litbang" ---
$d.Content.Find.Execute("llitbang", $true, $false, $false, $false, $false,
                         $true, 1, $false, "litbang", 2) | Out-Null

# --- Change 3: rewrite the "Jurnal ..." paragraph with its new wording,
#     split across four runs exactly as the target markup does. We use
#     InsertXML on the paragraph's range so the four <w:r> elements are
#     preserved distinctly (a plain .Text assignment would coalesce
#     same-formatted runs into a single run). ---
$target = "Jurnal di digunakan dikarenakan sesuai judulnya hanya saja berbeda metodenya. Dan mungkin parameter inputnya bisa digunakan sebagai warna buah itu sendiri."
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd() -eq $target) {
        $pStart = $para.Range.Start
        $pEnd = $para.Range.End
        $paraRange = $d.Range($pStart, $pEnd - 1)

        $xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:xml="http://www.w3.org/XML/1998/namespace">
<w:body>
<w:p>
<w:r><w:t xml:space="preserve">Jurnal </w:t></w:r>
<w:r><w:t xml:space="preserve">ini direkomendasikan oleh dosen pembimbing saya, selain itu mengapa saya menggunakan rujukan ini juga, itu dikarenakan </w:t></w:r>
<w:r><w:t xml:space="preserve">dosen pembimbing saya </w:t></w:r>
<w:r><w:t>bermitra dengan kementerian pertanian, maka dari itu. Penulis ini menggunakan ini</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
        $paraRange.InsertXML($xml)
        break
    }
}
